$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 2.1
$ws.Range("H5").Value = 3.4
$ws.Range("I5").Value = 3.4
$ws.Range("J5").Value = 2.75
$ws.Range("L5").Value = 3.75
$ws.Range("U5").Value = 1.67
$ws.Range("V5").Value = 2.1
$ws.Range("W5").Value = 8.5
$ws.Range("X5").Value = 11
$ws.Range("AA5").Value = 17
$ws.Range("AD5").Value = 6.5
$ws.Range("AI5").Value = 17
$ws.Range("AK5").Value = 34
$ws.Range("AL5").Value = 26
$ws.Range("AN5").Value = 4.33
$ws.Range("AZ5").Value = 51

# Row 6 updates
$ws.Range("G6").Value = 1.6
$ws.Range("H6").Value = 3.8
$ws.Range("I6").Value = 5.75
$ws.Range("J6").Value = 2.2
$ws.Range("L6").Value = 6
$ws.Range("N6").Value = 9.5
$ws.Range("Q6").Value = 2.03
$ws.Range("R6").Value = 1.83
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1.73
$ws.Range("Z6").Value = 11
$ws.Range("AE6").Value = 19
$ws.Range("AJ6").Value = 19
$ws.Range("AK6").Value = 67
$ws.Range("AQ6").Value = 26
$ws.Range("AU6").Value = 9
$ws.Range("AV6").Value = 67
$ws.Range("AY6").Value = 41
$ws.Range("AZ6").Value = 126
$ws.Range("BA6").Value = 151

# Row 8 updates
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 11
$ws.Range("AI8").Value = 26
